$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Part 1: Rework "Access Token" / "Client ID" / "493064913743651" /
# "Client Secret" paragraphs (the production credentials block):
#   - merge "Access Token" run fragments into a single run (keeps
#     its own paragraph's en-US language formatting)
#   - drop the en-US language formatting from "Client ID",
#     "493064913743651" and "Client Secret" paragraphs, and mark
#     "Client"/"Secret" as spell-check fragments via w:proofErr
# -----------------------------------------------------------------

$pAccessToken = $d.Paragraphs.Item(5)
$pShareHeading = $d.Paragraphs.Item(9)
$rng1 = $d.Range($pAccessToken.Range.Start, $pShareHeading.Range.Start)

$xml1 = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="2048">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Access Token APP_USR-493064913743651-062517-94d775dfe2c4f567f59dbedcb6e96646-160229129</w:t></w:r>
</w:p>
<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Client</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> ID</w:t></w:r></w:p>
<w:p><w:r><w:t>493064913743651</w:t></w:r></w:p>
<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Client</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Secret</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
$rng1.InsertXML($xml1)

# -----------------------------------------------------------------
# Part 2: Append a second, renewed set of credentials (sandbox /
# v2 style) plus the "share credentials" boilerplate at the end of
# the document, right before the last (empty) paragraph, which is
# kept as the very last paragraph but given an en-US language mark.
# -----------------------------------------------------------------

$pLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$rng2 = $pLast.Range
$rng2.Collapse(0)

$xml2 = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="2048">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:rPr><w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/><w:shd w:val="clear" w:color="auto" w:fill="EDEDED"/><w:lang w:val="en-US"/></w:rPr></w:pPr>
<w:r><w:rPr><w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/><w:shd w:val="clear" w:color="auto" w:fill="EDEDED"/><w:lang w:val="en-US"/></w:rPr><w:t>Public Key</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:rPr><w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/><w:shd w:val="clear" w:color="auto" w:fill="EDEDED"/><w:lang w:val="en-US"/></w:rPr></w:pPr>
<w:r><w:rPr><w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/><w:shd w:val="clear" w:color="auto" w:fill="EDEDED"/><w:lang w:val="en-US"/></w:rPr><w:t>APP_USR-2c6c878a-02c4-410f-89e0-cc2b244810aa</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Access Token</w:t></w:r>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>APP_USR-2496115398990356-062614-c613d17ba92fcfd69b4d88959f53a3b4-160229129</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Client ID</w:t></w:r>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>2496115398990356</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>2496115398990356</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Client Secret</w:t></w:r>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:t>JaeJrTrxg2ghwroxgQTsGJ286qVsGo4o</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>
<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Compartilhe as credenciais com um desenvolvedor</w:t></w:r>
</w:p>
<w:p>
<w:r><w:t>Se alguém está te ajudando a integrar os produtos do Mercado Pago, você pode compartilhar as credenciais da sua aplicação com essa pessoa de forma segura. Você pode desfazer essa ação excluindo a conta que recebeu o compartilhamento e renovando as credenciais.</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>
<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Compartilhar credenciais</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
$rng2.InsertXML($xml2)

# The trailing empty paragraph (kept intact by the insert above)
# picks up an en-US language mark in the final document.
$pTrailing = $d.Paragraphs.Item($d.Paragraphs.Count)
$pTrailing.Range.LanguageID = "en-US"
